$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the policy description labels in column A (rows 6-8). The shared-
# string table was reordered upstream so that "At least 0.7%..." and "Debt
# relief..." now come before "Expand Security Council...": reflect that by
# cycling the three labels accordingly.
$ws.Range("A6").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("A7").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("A8").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"

# Re-run autofit so the multi-line labels do not leave a stale explicit row
# height behind (keeps row metadata identical to a default render).
$ws.Rows(6).AutoFit()
$ws.Rows(7).AutoFit()
$ws.Rows(8).AutoFit()

# Update the numeric data (B2:N11) with the final rendered values from the
# refreshed "prepare & render" pipeline run.
$ws.Range("B2").Value = 0.809187541057965
$ws.Range("C2").Value = 0.843978206156072
$ws.Range("D2").Value = 0.866572127228228
$ws.Range("E2").Value = 0.827078685742989
$ws.Range("F2").Value = 0.888774505377363
$ws.Range("G2").Value = 0.788309742052631
$ws.Range("H2").Value = 0.813254042585487
$ws.Range("I2").Value = 0.853999575998437
$ws.Range("J2").Value = 0.773423212523683
$ws.Range("K2").Value = 0.806172053330117
$ws.Range("L2").Value = 0.801008544697843
$ws.Range("M2").Value = 0.859553623962162
$ws.Range("N2").Value = 0.767713744042883
$ws.Range("B3").Value = 0.793366965036613
$ws.Range("C3").Value = 0.820231998216116
$ws.Range("D3").Value = 0.811907205591862
$ws.Range("E3").Value = 0.806881985656633
$ws.Range("F3").Value = 0.883216508695824
$ws.Range("G3").Value = 0.716656390292452
$ws.Range("H3").Value = 0.80938575327201
$ws.Range("I3").Value = 0.845963888699056
$ws.Range("J3").Value = 0.752261706136369
$ws.Range("K3").Value = 0.805037944194611
$ws.Range("L3").Value = 0.829224529123518
$ws.Range("M3").Value = 0.870851551659494
$ws.Range("N3").Value = 0.740197166164522
$ws.Range("B4").Value = 0.748295871658231
$ws.Range("C4").Value = 0.748627231605022
$ws.Range("D4").Value = 0.719347570564758
$ws.Range("E4").Value = 0.733910266089133
$ws.Range("F4").Value = 0.83570773607823
$ws.Range("G4").Value = 0.717129396403951
$ws.Range("H4").Value = 0.770546743264546
$ws.Range("I4").Value = 0.724763201834602
$ws.Range("J4").Value = 0.667442126414247
$ws.Range("K4").Value = 0.727293734869913
$ws.Range("L4").Value = 0.866768469161838
$ws.Range("M4").Value = 0.894499688071821
$ws.Range("N4").Value = 0.697898274519676
$ws.Range("B5").Value = 0.699644803977167
$ws.Range("C5").Value = 0.732870289440397
$ws.Range("D5").Value = 0.776964615344406
$ws.Range("E5").Value = 0.696966943798898
$ws.Range("F5").Value = 0.779538293526919
$ws.Range("G5").Value = 0.612876561840658
$ws.Range("H5").Value = 0.738317533612799
$ws.Range("I5").Value = 0.747513876364388
$ws.Range("J5").Value = 0.718133301168214
$ws.Range("K5").Value = 0.58767644567848
$ws.Range("L5").Value = 0.72514029245591
$ws.Range("M5").Value = 0.814900578705803
$ws.Range("N5").Value = 0.674357186299684
$ws.Range("B6").Value = 0.698715666285492
$ws.Range("C6").Value = 0.689223510024568
$ws.Range("D6").Value = 0.663203246136222
$ws.Range("E6").Value = 0.668909891633291
$ws.Range("F6").Value = 0.787318727112064
$ws.Range("G6").Value = 0.586990880973224
$ws.Range("H6").Value = 0.772443771779816
$ws.Range("I6").Value = 0.654946097870933
$ws.Range("J6").Value = 0.641630244899629
$ws.Range("K6").Value = 0.622706835721545
$ws.Range("L6").Value = 0.825522087700358
$ws.Range("M6").Value = 0.863575793802146
$ws.Range("N6").Value = 0.665329447652658
$ws.Range("B7").Value = 0.696886644817168
$ws.Range("C7").Value = 0.698508510826339
$ws.Range("D7").Value = 0.642715315691589
$ws.Range("E7").Value = 0.600548950301388
$ws.Range("F7").Value = 0.812684805143714
$ws.Range("G7").Value = 0.794392466255568
$ws.Range("H7").Value = 0.71745735238834
$ws.Range("I7").Value = 0.724356614157177
$ws.Range("J7").Value = 0.650323736989535
$ws.Range("K7").Value = 0.683253184348732
$ws.Range("L7").Value = 0.746898541844184
$ws.Range("M7").Value = 0.883360143018035
$ws.Range("N7").Value = 0.666083774070859
$ws.Range("B8").Value = 0.694962460493881
$ws.Range("C8").Value = 0.761932435910776
$ws.Range("D8").Value = 0.720303245818442
$ws.Range("E8").Value = 0.756868069241599
$ws.Range("F8").Value = 0.802386607493625
$ws.Range("G8").Value = 0.7250772374002
$ws.Range("H8").Value = 0.763310673110362
$ws.Range("I8").Value = 0.78061912302082
$ws.Range("J8").Value = 0.721726707430624
$ws.Range("K8").Value = 0.67681938156
$ws.Range("L8").Value = 0.52836576507021
$ws.Range("M8").Value = 0.836791104476278
$ws.Range("N8").Value = 0.671764356435874
$ws.Range("B9").Value = 0.683174203642518
$ws.Range("C9").Value = 0.69192411456902
$ws.Range("D9").Value = 0.677768220232769
$ws.Range("E9").Value = 0.688283637483319
$ws.Range("F9").Value = 0.757717193334531
$ws.Range("G9").Value = 0.632201156950068
$ws.Range("H9").Value = 0.73213923408807
$ws.Range("I9").Value = 0.668233479894725
$ws.Range("J9").Value = 0.645323472858561
$ws.Range("K9").Value = 0.59008286280892
$ws.Range("L9").Value = 0.875589503017251
$ws.Range("M9").Value = 0.858825470219245
$ws.Range("N9").Value = 0.614162213651268
$ws.Range("B10").Value = 0.682119272063001
$ws.Range("C10").Value = 0.746931207916325
$ws.Range("D10").Value = 0.749211487143156
$ws.Range("E10").Value = 0.739300715467795
$ws.Range("F10").Value = 0.848419375204591
$ws.Range("G10").Value = 0.659131063582329
$ws.Range("H10").Value = 0.695127487350912
$ws.Range("I10").Value = 0.739704530251875
$ws.Range("J10").Value = 0.625377265896113
$ws.Range("K10").Value = 0.7245399038347
$ws.Range("L10").Value = 0.50012458723726
$ws.Range("M10").Value = 0.770097780765828
$ws.Range("N10").Value = 0.655816868493662
$ws.Range("B11").Value = 0.526319576644156
$ws.Range("C11").Value = 0.551475483086733
$ws.Range("D11").Value = 0.616758276280248
$ws.Range("E11").Value = 0.542904806999543
$ws.Range("F11").Value = 0.555847318045528
$ws.Range("G11").Value = 0.527788999798779
$ws.Range("H11").Value = 0.538544964808905
$ws.Range("I11").Value = 0.535237211156328
$ws.Range("J11").Value = 0.506594819956859
$ws.Range("K11").Value = 0.46350801136536
$ws.Range("L11").Value = 0.514246894490421
$ws.Range("M11").Value = 0.695189510992556
$ws.Range("N11").Value = 0.507370211647722
